$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 12): Date, Completed Chapters, and the
# shared formula continues into C12.
$ws.Range("A12").Value = 45795
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat

$ws.Range("B12").Value = 65

$ws.Range("C12").Formula = "=B12/200*100"

# Move the active cell selection to G16 to match the saved state.
$ws.Range("G16").Select()
